# Add four new "Statement of the Problem" list items after the existing
# "What would help regulate the products sold in school?" bullet. Each new
# paragraph is appended as a sibling list item so it inherits the
# ListParagraph style / numbering (numId 1) / spacing / run formatting
# (Times New Roman, 12pt) already used by the preceding bullets.

$d = $word.ActiveDocument

$newBullets = @(
    "How can the students maximize their break or lunch time in terms of purchasing food?",
    "How can the vendors accommodate more customers at a time?",
    "What can be a more convenient way to do transactions in terms of payment?",
    "What can help with the inventory process which the vendors do at the end of the day?"
)

foreach ($bulletText in $newBullets) {
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $insertionRange = $lastPara.Range
    $insertionRange.Collapse(0)
    $insertionRange.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newRange = $newPara.Range
    $newRange.Collapse(0)
    $newRange.Text = $bulletText
}
